# Auto-generated Excel COM-interop script
# Applies market-price / profit data refresh to the Brynhildr_Profits workbook
# across all 8 sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
#
# For each affected row, cells are either updated to a new numeric value,
# or cleared entirely (when the source diff shows the cell disappearing),
# or newly populated (when the source diff shows the cell appearing).

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 104.6
$ws.Range("I2").Value = 128
$ws.Range("K2").Value = 128
$ws.Range("M2").Value = -15
$ws.Range("H28").Value = 446.6207
$ws.Range("I28").Value = 462.68182
$ws.Range("K28").Value = 462.68182
$ws.Range("M28").Value = 22.31817999999998
$ws.Range("H32").Value = 2000
$ws.Range("I32").Value = 2000
$ws.Range("J32").Value = 2000
$ws.Range("K32").Value = 2000
$ws.Range("L32").Value = 2000
$ws.Range("M32").Value = -1674
$ws.Range("N32").Value = -2652
$ws.Range("H96").Value = 1644.0526
$ws.Range("I96").Value = 402.25
$ws.Range("J96").Value = 2547.182
$ws.Range("K96").Value = 1206.75
$ws.Range("L96").Value = 7641.545999999999
$ws.Range("M96").Value = 166.25
$ws.Range("N96").Value = -10387.546
$ws.Range("H99").Value = 190
$ws.Range("I99").Value = 190
$ws.Range("K99").Value = 570
$ws.Range("M99").Value = 928
$ws.Range("H100").Value = 3056.162
$ws.Range("I100").Value = 936.45
$ws.Range("K100").Value = 936.45
$ws.Range("M100").Value = -395.45
$ws.Range("H116").Value = 39405.445
$ws.Range("I116").Value = 29916.666
$ws.Range("K116").Value = 29916.666
$ws.Range("M116").Value = -26474.666
$ws.Range("H132").Value = 20558.273
$ws.Range("I132").Value = 27566.375
$ws.Range("J132").Value = 1870
$ws.Range("K132").Value = 82699.125
$ws.Range("L132").Value = 5610
$ws.Range("M132").Value = -80169.125
$ws.Range("N132").Value = -10670
$ws.Range("H137").Value = 3589
$ws.Range("I137").Value = 1593.3556
$ws.Range("K137").Value = 4780.066800000001
$ws.Range("M137").Value = -2230.066800000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 813.56525
$ws.Range("I2").Value = 839.6111
$ws.Range("K2").Value = 839.6111
$ws.Range("M2").Value = -726.6111
$ws.Range("H61").Value = 3725.6155
$ws.Range("I61").Value = 3727.75
$ws.Range("K61").Value = 3727.75
$ws.Range("M61").Value = -3515.75
$ws.Range("H116").Value = 813.56525
$ws.Range("I116").Value = 839.6111
$ws.Range("K116").Value = 839.6111
$ws.Range("M116").Value = 1454.3889
$ws.Range("H136").Value = 3725.6155
$ws.Range("I136").Value = 3727.75
$ws.Range("K136").Value = 11183.25
$ws.Range("M136").Value = -8633.25

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 813.56525
$ws.Range("I3").Value = 839.6111
$ws.Range("K3").Value = 839.6111
$ws.Range("M3").Value = -725.6111
$ws.Range("H60").Value = 87500
$ws.Range("J60").Value = 87500
$ws.Range("L60").Value = 87500
$ws.Range("N60").Value = -88698
$ws.Range("H94").Value = 5087.8945
$ws.Range("I94").Value = 4416.875
$ws.Range("K94").Value = 4416.875
$ws.Range("M94").Value = -3965.875
$ws.Range("H99").Value = 34933.332
$ws.Range("I99").Value = 34933.332
$ws.Range("K99").Value = 34933.332
$ws.Range("M99").Value = -33435.332
$ws.Range("H105").Value = 7956
$ws.Range("I105").Value = 1949.5
$ws.Range("K105").Value = 1949.5
$ws.Range("M105").Value = -202.5
$ws.Range("H134").Value = 10066.833
$ws.Range("I134").Value = 10066.833
$ws.Range("K134").Value = 30200.499
$ws.Range("M134").Value = -27665.499
$ws.Range("H135").Value = 66691.39999999999
$ws.Range("J135").Value = 66691.39999999999
$ws.Range("L135").Value = 66691.39999999999
$ws.Range("N135").Value = -76831.39999999999
$ws.Range("H137").Value = 0
$ws.Range("I137").Value = 0
$ws.Range("J137").Value = 0
$ws.Range("K137").Value = 0
$ws.Range("L137").Value = 0
$ws.Range("M137").ClearContents()
$ws.Range("N137").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H21").Value = 0
$ws.Range("I21").Value = 0
$ws.Range("K21").Value = 0
$ws.Range("M21").ClearContents()
$ws.Range("H25").Value = 9174
$ws.Range("I25").Value = 9174
$ws.Range("J25").Value = 0
$ws.Range("K25").Value = 9174
$ws.Range("L25").Value = 0
$ws.Range("M25").Value = -9000
$ws.Range("N25").ClearContents()
$ws.Range("H99").Value = 22486.9
$ws.Range("J99").Value = 2838
$ws.Range("L99").Value = 2838
$ws.Range("N99").Value = -5834
$ws.Range("H107").Value = 1185.5
$ws.Range("I107").Value = 1055.9333
$ws.Range("K107").Value = 1055.9333
$ws.Range("M107").Value = 864.0667000000001
$ws.Range("H126").Value = 22486.9
$ws.Range("J126").Value = 2838
$ws.Range("L126").Value = 8514
$ws.Range("N126").Value = -13454

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H9").Value = 1555125.1
$ws.Range("I9").Value = 1093860
$ws.Range("J9").Value = 1801133.2
$ws.Range("K9").Value = 3281580
$ws.Range("L9").Value = 5403399.6
$ws.Range("M9").Value = -3281356
$ws.Range("N9").Value = -5403847.6
$ws.Range("H12").Value = 582.5
$ws.Range("I12").Value = 312.75
$ws.Range("J12").Value = 852.25
$ws.Range("K12").Value = 938.25
$ws.Range("L12").Value = 2556.75
$ws.Range("M12").Value = -765.25
$ws.Range("N12").Value = -2902.75
$ws.Range("H37").Value = 46665
$ws.Range("J37").Value = 46665
$ws.Range("L37").Value = 139995
$ws.Range("N37").Value = -140219
$ws.Range("H113").Value = 24464.52
$ws.Range("J113").Value = 35785.59
$ws.Range("L113").Value = 107356.77
$ws.Range("N113").Value = -111696.77
$ws.Range("H139").Value = 9731.929
$ws.Range("I139").Value = 6292
$ws.Range("J139").Value = 12311.875
$ws.Range("K139").Value = 18876
$ws.Range("L139").Value = 36935.625
$ws.Range("M139").Value = -13736
$ws.Range("N139").Value = -47215.625

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H21").Value = 802756.2
$ws.Range("I21").Value = 802756.2
$ws.Range("K21").Value = 802756.2
$ws.Range("M21").Value = -802583.2
$ws.Range("H30").Value = 802756.2
$ws.Range("I30").Value = 802756.2
$ws.Range("K30").Value = 802756.2
$ws.Range("M30").Value = -802651.2
$ws.Range("H132").Value = 18135.6
$ws.Range("I132").Value = 18135.6
$ws.Range("K132").Value = 54406.8
$ws.Range("M132").Value = -51876.8

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1367.6333
$ws.Range("I16").Value = 1465.32
$ws.Range("J16").Value = 879.2
$ws.Range("K16").Value = 1465.32
$ws.Range("L16").Value = 879.2
$ws.Range("M16").Value = -1295.32
$ws.Range("N16").Value = -1219.2
$ws.Range("H23").Value = 25000
$ws.Range("I23").Value = 25000
$ws.Range("J23").Value = 0
$ws.Range("K23").Value = 25000
$ws.Range("L23").Value = 0
$ws.Range("M23").Value = -24770
$ws.Range("N23").ClearContents()
$ws.Range("H93").Value = 4017.375
$ws.Range("I93").Value = 1632.4
$ws.Range("K93").Value = 1632.4
$ws.Range("M93").Value = -384.4000000000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H30").Value = 25500
$ws.Range("I30").Value = 9000
$ws.Range("K30").Value = 9000
$ws.Range("M30").Value = -8893
$ws.Range("H80").Value = 27500
$ws.Range("I80").Value = 5000
$ws.Range("J80").Value = 50000
$ws.Range("K80").Value = 5000
$ws.Range("L80").Value = 50000
$ws.Range("M80").Value = -4002
$ws.Range("N80").Value = -51996
$ws.Range("H81").Value = 8335.5
$ws.Range("I81").Value = 2646.4285
$ws.Range("K81").Value = 5292.857
$ws.Range("M81").Value = -4231.857
$ws.Range("H82").Value = 55000
$ws.Range("J82").Value = 55000
$ws.Range("L82").Value = 55000
$ws.Range("N82").Value = -55766
$ws.Range("H83").Value = 27500
$ws.Range("I83").Value = 5000
$ws.Range("J83").Value = 50000
$ws.Range("K83").Value = 15000
$ws.Range("L83").Value = 150000
$ws.Range("M83").Value = -10008
$ws.Range("N83").Value = -159984
$ws.Range("H84").Value = 8335.5
$ws.Range("I84").Value = 2646.4285
$ws.Range("K84").Value = 26464.285
$ws.Range("M84").Value = -21160.285
$ws.Range("H85").Value = 55000
$ws.Range("J85").Value = 55000
$ws.Range("L85").Value = 55000
$ws.Range("N85").Value = -57652
$ws.Range("H136").Value = 3476.8064
$ws.Range("J136").Value = 684.5
$ws.Range("L136").Value = 2053.5
$ws.Range("N136").Value = -7153.5
